$wb = $excel.ActiveWorkbook

# --- "loads" sheet: rebuild the table with two new leading columns
# (v_nom_kv, s_base_mva) after "name" and two new trailing columns
# (g_shunt_pu, b_shunt_pu) after "bus_idx".
$loads = $wb.Worksheets.Item("loads")

$loads.Cells.Item(1, 1).Value = "name"
$loads.Cells.Item(1, 2).Value = "v_nom_kv"
$loads.Cells.Item(1, 3).Value = "s_base_mva"
$loads.Cells.Item(1, 4).Value = "v_nom_pu"
$loads.Cells.Item(1, 5).Value = "p_nom_mw"
$loads.Cells.Item(1, 6).Value = "q_nom_mvar"
$loads.Cells.Item(1, 7).Value = "bus_idx"
$loads.Cells.Item(1, 8).Value = "g_shunt_pu"
$loads.Cells.Item(1, 9).Value = "b_shunt_pu"

$loads.Cells.Item(2, 1).Value = "Load 1"
$loads.Cells.Item(2, 2).Value = 132
$loads.Cells.Item(2, 3).Value = 100
$loads.Cells.Item(2, 4).Value = 1
$loads.Cells.Item(2, 5).Value = 20
$loads.Cells.Item(2, 6).Value = 10
$loads.Cells.Item(2, 7).Value = 1
$loads.Cells.Item(2, 8).Value = 0
$loads.Cells.Item(2, 9).Value = 0

$loads.Cells.Item(3, 1).Value = "Load 2"
$loads.Cells.Item(3, 2).Value = 22
$loads.Cells.Item(3, 3).Value = 100
$loads.Cells.Item(3, 4).Value = 1
$loads.Cells.Item(3, 5).Value = 20
$loads.Cells.Item(3, 6).Value = 10
$loads.Cells.Item(3, 7).Value = 4
$loads.Cells.Item(3, 8).Value = 0
$loads.Cells.Item(3, 9).Value = 0

# Make "loads" the active / tab-selected sheet, with the new trailing cell selected
$loads.Select() | Out-Null
$loads.Range("I4").Select() | Out-Null
